$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("H 72") - all rows below shift up by one
$ws.Rows.Item(2).Delete()

# Apply the updated "missing data" mask / values for columns C, D, F
$ws.Range("F2").ClearContents()
$ws.Range("C3").Value = 11.7
$ws.Range("D3").Value = -13.5
$ws.Range("C4").ClearContents()
$ws.Range("C5").Value = 11.2
$ws.Range("D5").ClearContents()
$ws.Range("F5").Value = 0.71055
$ws.Range("C8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("D10").Value = -13.8
$ws.Range("F10").Value = 0.7105
$ws.Range("D11").Value = -13.9
$ws.Range("D12").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("F13").ClearContents()
$ws.Range("D14").Value = -15.5
$ws.Range("C15").Value = 12
$ws.Range("C16").ClearContents()
$ws.Range("D16").ClearContents()
$ws.Range("F18").Value = 0.71073
$ws.Range("C19").Value = 12.5
$ws.Range("F19").Value = 0.7107599999999999
$ws.Range("F20").Value = 0.7106
$ws.Range("C21").ClearContents()
$ws.Range("F23").ClearContents()
$ws.Range("F24").ClearContents()
$ws.Range("F25").ClearContents()
$ws.Range("F26").Value = 0.70925
$ws.Range("C27").Value = 13.5
$ws.Range("F27").Value = 0.7092000000000001
$ws.Range("D29").Value = -13.9
$ws.Range("F29").ClearContents()
$ws.Range("C30").ClearContents()
$ws.Range("F30").ClearContents()
$ws.Range("C31").Value = 10.7
$ws.Range("D31").ClearContents()
$ws.Range("C32").ClearContents()
$ws.Range("F33").Value = 0.72961
$ws.Range("F36").ClearContents()
$ws.Range("D37").Value = -14.3
$ws.Range("C39").Value = 11.7
$ws.Range("D39").ClearContents()
$ws.Range("C40").ClearContents()
$ws.Range("F41").Value = 0.7111499999999999
$ws.Range("F42").Value = 0.7111499999999999
$ws.Range("F43").Value = 0.71152
$ws.Range("F44").ClearContents()
$ws.Range("C45").Value = 11.7
$ws.Range("D45").Value = -14.9
$ws.Range("F45").ClearContents()
$ws.Range("C46").ClearContents()
$ws.Range("F46").ClearContents()
$ws.Range("D47").ClearContents()
$ws.Range("D49").Value = -13.7
$ws.Range("C50").Value = 10.7
$ws.Range("D50").Value = -15.6
$ws.Range("D51").ClearContents()
$ws.Range("C52").ClearContents()
$ws.Range("D52").ClearContents()
$ws.Range("D53").Value = -12.9
$ws.Range("D55").ClearContents()
$ws.Range("F55").Value = 0.71492
$ws.Range("C56").Value = 11.9
$ws.Range("C57").ClearContents()
$ws.Range("C58").Value = 11.2
$ws.Range("D58").Value = -13
$ws.Range("F58").ClearContents()
$ws.Range("C59").ClearContents()
$ws.Range("D59").Value = -13.6
$ws.Range("F60").Value = 0.70948
$ws.Range("C61").Value = 10.5
$ws.Range("C62").ClearContents()
$ws.Range("D62").ClearContents()
